# Insert one new data row into the "Vega Modelo de Temuco - Melón" sheet.
#
# The canonical diff shows a brand-new record inserted at row 760, which
# pushes every following row (760-830) down by one (761-831) and grows the
# used range from A1:R830 to A1:R831. We reproduce that by inserting a row
# at position 760 (which shifts everything below it down, copying the
# formatting of the surrounding rows) and then writing the new record's
# values into the freshly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 760 - this shifts rows 760:830
# down to 761:831 and keeps their values/formatting intact.
$ws.Rows.Item(760).Insert()

# Populate the newly inserted row 760 with the new record.
$ws.Cells.Item(760, 1).Value  = 10
$ws.Cells.Item(760, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(760, 3).Value  = "La Araucanía"
$ws.Cells.Item(760, 4).Value  = 45041
$ws.Cells.Item(760, 5).Value  = 9
$ws.Cells.Item(760, 6).Value  = 100112027
$ws.Cells.Item(760, 7).Value  = "Melón"
$ws.Cells.Item(760, 8).Value  = "Tuna"
$ws.Cells.Item(760, 9).Value  = "Primera"
$ws.Cells.Item(760, 10).Value = 40
$ws.Cells.Item(760, 11).Value = 25000
$ws.Cells.Item(760, 12).Value = 25000
$ws.Cells.Item(760, 13).Value = 25000
$ws.Cells.Item(760, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(760, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(760, 16).Value = 1389
$ws.Cells.Item(760, 17).Value = 18
$ws.Cells.Item(760, 18).Value = "Hortaliza"
